# Update column C ("Förändrad") date value from 2023-10-04 (45203) to
# 2023-10-05 (45204) for every data row (rows 2 through the last used row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 505 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45204
}
